$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# Update the reserve / station-key bounding box coordinates after converting all
# shapefiles to WGS 84 (EPSG 4269).
$ws.Range("A2").Value = -122.5601
$ws.Range("B2").Value = -122.5325

$ws.Range("A3").Value = 37.9392
$ws.Range("B3").Value = 37.9606

$ws.Range("A4").Value = -121.9468
$ws.Range("B4").Value = -121.9744

$ws.Range("A5").Value = 38.2714
$ws.Range("B5").Value = 38.25

# Touch the formatting on B6 (left blank, no bounding box value there).
$ws.Range("B6").ClearFormats()
